# Applies the diff: inserts a new top headline row (PFAS story) into the
# rotating 9-row news list (rows 3-11), shifting existing rows down by one
# and dropping the previous last row. Also refreshes the shared thumbnail
# image id from 345866 to 345911 for all list rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "New rules are coming for PFAS in drinking water. See how your water tests."
    "With Bowman’s loss to Latimer, the NY Democratic establishment strikes back"
    "Bowman loses to Latimer in NY Democratic primary shaped by Israel-Hamas war"
    "Admission for Black students at NYC’s specialized high schools ticks up slightly, but remains low"
    "NYPD inspector accused of sloppy coverup of DUI crash after boozy night at Midtown bar"
    "NYCHA’s not just a landlord. It's going to train young adults to work in construction."
    "James Beard Foundation announces Best Chef in New York"
    "NYC public schools will educate students and their families about safe gun storage"
    "How to start running in NYC: Fun routes, safety tips and where to get proper shoes"
)

$descriptions = @(
    "A Gothamist analysis shows that the drinking water of 4.3 million New York and New Jersey residents has tested positive for so-called “forever chemicals.” See if you’re one of them."
    "Westchester’s George Latimer secured a quick and decisive win over Rep. Jamaal Bowman Tuesday night. Here’s what the outcome tells us."
    "The 16th Congressional District that spans from northern Bronx to Westchester now ranks as the most expensive House primary in American history."
    "“Instead of pushing for systemwide change, this administration applauds a mere 1% increase in offers to Black and Latinx students,” one advocate says."
    "Prosecutors said Deputy Inspector Paul Zangrilli scrambled to destroy evidence after he let his drunk girlfriend drive his unmarked police car."
    "The `$1.3 million grant from the U.S. Labor Department will fund the workforce development program for 40 months."
    "New York was iced out of many major categories this year."
    "The announcement comes after a 14-year-old was accidentally shot and killed by his cousin, who was playing with a gun."
    "Go extremely slow, try to find a car-free area — even if it’s just a little park to run laps around — and don’t be afraid to dress like a highlighter."
)

$images = @(
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
    "https://api-prod.gothamist.com/images/345911/fill-318x212%7Cformat-webp%7Cwebpquality-80/"
)

$counts = @(1, 0, 0, 0, 0, 0, 0, 0, 0)
$containsMoney = @($false, $false, $false, $false, $false, $true, $false, $false, $false)

for ($i = 0; $i -lt 9; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $titles[$i]
    $ws.Cells.Item($row, 3).Value = $descriptions[$i]
    $ws.Cells.Item($row, 4).Value = $images[$i]
    $ws.Cells.Item($row, 5).Value = $counts[$i]
    $ws.Cells.Item($row, 6).Value = $containsMoney[$i]
}
